$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the task description ("ist" -> "its") for the "Implement the virus..." task (row 3, col A)
$ws.Range("A3").Value = "Implement the virus with all its math and physiks"

# Row 2: Curr. Est. (C2) becomes a fixed value instead of a formula referencing B2
$ws.Range("C2").Value = 8

# Row 2: Effort (D2) updated to 5.5 with a one-decimal number format
$ws.Range("D2").NumberFormat = "0.0"
$ws.Range("D2").Value = 5.5

# Row 2: Remain (E2) keeps its formula (C2-D2) but adopts the same one-decimal number format
$ws.Range("E2").NumberFormat = "0.0"

# Update the last-saved selection/active cell to match the edited workbook
$ws.Range("D15").Select()
